
# --------------------------------------------------------------------------
# 2021-05 Victorian Outbreak Paths.xlsx
#
# Adds one new outbreak-path record (2021-06-26, Southbank "S 12") to
# Table1 on "Sheet1", and extends the Date_Colours gradient table on
# "Date Colours" with two more dates (2021-06-25 and 2021-06-26),
# which shifts the interpolated "Colour Code" gradient for every
# existing date row as well.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet1 / Table1 — append the new outbreak-path row
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$table1 = $ws1.ListObjects.Item("Table1")

$newPathRow = $table1.ListRows.Add()
$r = $newPathRow.Range.Row

$ws1.Cells.Item($r, 1).Value = 44373
$ws1.Cells.Item($r, 1).NumberFormat = "d-mmm"
$ws1.Cells.Item($r, 2).Value = "A / S"
$ws1.Cells.Item($r, 3).Value = "S 12"
$ws1.Cells.Item($r, 4).Value = "Southbank"
$ws1.Cells.Item($r, 6).Value = "Southbank Apartments"
$ws1.Cells.Item($r, 7).Value = "Kappa (B.1.617.1)"

# Match the author's final selection/cursor position
$ws1.Activate()
$ws1.Range("D109").Select()

# ---------------------------------------------------------------------
# 2) "Date Colours" / Date_Colours table — refresh the gradient and
#    append the two new dates
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Date Colours")
$table2 = $ws2.ListObjects.Item("Date_Colours")

# The "Colour Code" gradient (column B) is recomputed across all 33
# dates now that the table spans 2021-05-25 .. 2021-06-26 (was
# 2021-05-25 .. 2021-06-24). Column C ("Colour Code dbRaevn") and D
# ("Shape dbRaevn") are unaffected for the existing rows.
$gradient = @(
    "#fbfaff",
    "#f6f5ff",
    "#f2f0fe",
    "#edebfe",
    "#e9e6fe",
    "#e4e1fe",
    "#e0dcfd",
    "#dbd8fd",
    "#d7d3fd",
    "#d2cefc",
    "#cdc9fc",
    "#c9c4fc",
    "#c4c0fb",
    "#bfbbfb",
    "#bab6fa",
    "#b5b1fa",
    "#b0adfa",
    "#aba8f9",
    "#a6a4f9",
    "#a19ff8",
    "#9b9af8",
    "#9696f7",
    "#9091f7",
    "#8a8df6",
    "#8588f6",
    "#7e84f5",
    "#7880f5",
    "#717bf4",
    "#6a77f3",
    "#6373f3",
    "#5b6ef2"
)

for ($i = 0; $i -lt $gradient.Count; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 2).Value = $gradient[$i]
}

# Two new date rows appended to the table
$newColourRow1 = $table2.ListRows.Add()
$rc1 = $newColourRow1.Range.Row
$ws2.Cells.Item($rc1, 1).Value = 44372
$ws2.Cells.Item($rc1, 1).NumberFormat = "d-mmm"
$ws2.Cells.Item($rc1, 2).Value = "#536af2"
$ws2.Cells.Item($rc1, 4).Value = "Diamond"

$newColourRow2 = $table2.ListRows.Add()
$rc2 = $newColourRow2.Range.Row
$ws2.Cells.Item($rc2, 1).Value = 44373
$ws2.Cells.Item($rc2, 1).NumberFormat = "d-mmm"
$ws2.Cells.Item($rc2, 2).Value = "#4966f1"
$ws2.Cells.Item($rc2, 3).Value = "#FFCFAF"
$ws2.Cells.Item($rc2, 4).Value = "Diamond"
